# Update cryptos list values per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.651.99'
$ws.Range("E2").Value = '  +1.49%  '

$ws.Range("D3").Value = '1.602.51'
$ws.Range("E3").Value = '  +1.27%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = "'212.46"
$ws.Range("E5").Value = '  -0.14%  '

$ws.Range("E6").Value = '  +0.96%  '

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").Value = "'27.92"
$ws.Range("E8").Value = '  +5.78%  '

$ws.Range("E9").Value = '  +1.55%  '

$ws.Range("D10").Value = "'0.0602"
$ws.Range("E10").Value = '  +1.45%  '

$ws.Range("D11").Value = "'0.0909"
$ws.Range("E11").Value = '  +0.45%  '

$ws.Range("D12").Value = '1.832.53'
$ws.Range("E12").Value = '  +1.36%  '

$ws.Range("D13").Value = '1.597.23'
$ws.Range("E13").Value = '  +1.16%  '

$ws.Range("D14").Value = "'0.547"
$ws.Range("E14").Value = '  +4.22%  '

$ws.Range("D15").Value = '29.679.08'
$ws.Range("E15").Value = '  +1.53%  '

$ws.Range("E16").Value = '  +0.83%  '

$ws.Range("D17").Value = "'63.95"
$ws.Range("E17").Value = '  +1.86%  '

$ws.Range("D18").Value = "'242.57"
$ws.Range("E18").Value = '  +1.73%  '

$ws.Range("D19").Value = "'7.79"
$ws.Range("E19").Value = '  +4.13%  '

$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("E22").Value = '  +0.82%  '

$ws.Range("E23").Value = '  +2.69%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("D25").Value = "'155.42"
$ws.Range("E25").Value = '  +0.64%  '

$ws.Range("D26").Value = "'15.49"
$ws.Range("E26").Value = '  +2.11%  '

$ws.Range("D27").Value = "'0.109"
$ws.Range("E27").Value = '  +0.43%  '

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("E30").Value = '  +2.51%  '

$ws.Range("E31").Value = '  +0.16%  '

$ws.Range("E32").Value = '  +0.55%  '

$ws.Range("E33").Value = '  +3.52%  '

$ws.Range("D34").Value = '1.425.89'
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("E35").Value = '  +3.92%  '

$ws.Range("E36").Value = '  +3.54%  '

$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("D38").Value = "'2.29"
$ws.Range("E38").Value = '  -0.49%  '

$ws.Range("D40").Value = "'58.45"
$ws.Range("E40").Value = '  +8.22%  '

$ws.Range("D41").Value = "'0.547"
$ws.Range("E41").Value = '  +2.81%  '

$ws.Range("D42").Value = "'0.0498"
$ws.Range("E42").Value = '  +5.95%  '

$ws.Range("E43").Value = '  +0.49%  '

$ws.Range("E44").Value = '  +2.98%  '

$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("D46").Value = "'66.41"
$ws.Range("E46").Value = '  +2.85%  '

$ws.Range("D47").Value = "'0.972"
$ws.Range("E47").Value = '  +15.72%  '

$ws.Range("D48").Value = "'5.34"
$ws.Range("E48").Value = '  -0.08%  '

$ws.Range("D49").Value = '1.743.43'
$ws.Range("E49").Value = '  +1.35%  '

$ws.Range("D50").Value = "'86.83"
$ws.Range("E50").Value = '  +1.54%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.0524"
$ws.Range("E51").Value = '  +1.79%  '

# Remove the quote-prefix style Excel applies when forcing text via a
# leading apostrophe, so D2:D51's style stays identical to before.
$ws.Range("D2:D51").ClearFormats()
